# Auto update Excel log
# Appends newly captured sensor readings (2026-01-30, ~12:54-12:55) to each log sheet,
# forcing text formatting so values round-trip exactly as logged (no date/number coercion).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALERTS")
$rng = $ws.Range("A3:F3")
$rng.NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "2026-01-30"
$ws.Cells.Item(3, 2).Value = "12:55:19"
$ws.Cells.Item(3, 3).Value = "12:00"
$ws.Cells.Item(3, 4).Value = "Bathroom"
$ws.Cells.Item(3, 5).Value = "MODERATE"
$ws.Cells.Item(3, 6).Value = "MODERATE ALERT: Bathroom occupied, no motion > 40s."

$ws = $wb.Worksheets.Item("PIR")
$rng = $ws.Range("A31:F31")
$rng.NumberFormat = "@"
$ws.Cells.Item(31, 1).Value = "2026-01-30"
$ws.Cells.Item(31, 2).Value = "12:54:21"
$ws.Cells.Item(31, 3).Value = "12:00"
$ws.Cells.Item(31, 4).Value = "Bathroom"
$ws.Cells.Item(31, 5).Value = "No Motion"
$ws.Cells.Item(31, 6).Value = "Inactive"
$rng = $ws.Range("A32:F32")
$rng.NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "2026-01-30"
$ws.Cells.Item(32, 2).Value = "12:54:21"
$ws.Cells.Item(32, 3).Value = "12:00"
$ws.Cells.Item(32, 4).Value = "Bathroom"
$ws.Cells.Item(32, 5).Value = "No Motion"
$ws.Cells.Item(32, 6).Value = "Inactive"
$rng = $ws.Range("A33:F33")
$rng.NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = "2026-01-30"
$ws.Cells.Item(33, 2).Value = "12:54:24"
$ws.Cells.Item(33, 3).Value = "12:00"
$ws.Cells.Item(33, 4).Value = "Bathroom"
$ws.Cells.Item(33, 5).Value = "No Motion"
$ws.Cells.Item(33, 6).Value = "Inactive"
$rng = $ws.Range("A34:F34")
$rng.NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = "2026-01-30"
$ws.Cells.Item(34, 2).Value = "12:54:29"
$ws.Cells.Item(34, 3).Value = "12:00"
$ws.Cells.Item(34, 4).Value = "Bathroom"
$ws.Cells.Item(34, 5).Value = "No Motion"
$ws.Cells.Item(34, 6).Value = "Inactive"
$rng = $ws.Range("A35:F35")
$rng.NumberFormat = "@"
$ws.Cells.Item(35, 1).Value = "2026-01-30"
$ws.Cells.Item(35, 2).Value = "12:54:33"
$ws.Cells.Item(35, 3).Value = "12:00"
$ws.Cells.Item(35, 4).Value = "Living Room"
$ws.Cells.Item(35, 5).Value = "RECOVERY_DETECTION"
$ws.Cells.Item(35, 6).Value = "Inactive"
$rng = $ws.Range("A36:F36")
$rng.NumberFormat = "@"
$ws.Cells.Item(36, 1).Value = "2026-01-30"
$ws.Cells.Item(36, 2).Value = "12:54:35"
$ws.Cells.Item(36, 3).Value = "12:00"
$ws.Cells.Item(36, 4).Value = "Bathroom"
$ws.Cells.Item(36, 5).Value = "No Motion"
$ws.Cells.Item(36, 6).Value = "Inactive"
$rng = $ws.Range("A37:F37")
$rng.NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = "2026-01-30"
$ws.Cells.Item(37, 2).Value = "12:54:40"
$ws.Cells.Item(37, 3).Value = "12:00"
$ws.Cells.Item(37, 4).Value = "Bathroom"
$ws.Cells.Item(37, 5).Value = "No Motion"
$ws.Cells.Item(37, 6).Value = "Inactive"
$rng = $ws.Range("A38:F38")
$rng.NumberFormat = "@"
$ws.Cells.Item(38, 1).Value = "2026-01-30"
$ws.Cells.Item(38, 2).Value = "12:54:45"
$ws.Cells.Item(38, 3).Value = "12:00"
$ws.Cells.Item(38, 4).Value = "Bathroom"
$ws.Cells.Item(38, 5).Value = "No Motion"
$ws.Cells.Item(38, 6).Value = "Inactive"
$rng = $ws.Range("A39:F39")
$rng.NumberFormat = "@"
$ws.Cells.Item(39, 1).Value = "2026-01-30"
$ws.Cells.Item(39, 2).Value = "12:54:50"
$ws.Cells.Item(39, 3).Value = "12:00"
$ws.Cells.Item(39, 4).Value = "Bathroom"
$ws.Cells.Item(39, 5).Value = "No Motion"
$ws.Cells.Item(39, 6).Value = "Inactive"
$rng = $ws.Range("A40:F40")
$rng.NumberFormat = "@"
$ws.Cells.Item(40, 1).Value = "2026-01-30"
$ws.Cells.Item(40, 2).Value = "12:54:55"
$ws.Cells.Item(40, 3).Value = "12:00"
$ws.Cells.Item(40, 4).Value = "Bathroom"
$ws.Cells.Item(40, 5).Value = "No Motion"
$ws.Cells.Item(40, 6).Value = "Inactive"
$rng = $ws.Range("A41:F41")
$rng.NumberFormat = "@"
$ws.Cells.Item(41, 1).Value = "2026-01-30"
$ws.Cells.Item(41, 2).Value = "12:55:00"
$ws.Cells.Item(41, 3).Value = "12:00"
$ws.Cells.Item(41, 4).Value = "Bathroom"
$ws.Cells.Item(41, 5).Value = "No Motion"
$ws.Cells.Item(41, 6).Value = "Inactive"
$rng = $ws.Range("A42:F42")
$rng.NumberFormat = "@"
$ws.Cells.Item(42, 1).Value = "2026-01-30"
$ws.Cells.Item(42, 2).Value = "12:55:05"
$ws.Cells.Item(42, 3).Value = "12:00"
$ws.Cells.Item(42, 4).Value = "Bathroom"
$ws.Cells.Item(42, 5).Value = "No Motion"
$ws.Cells.Item(42, 6).Value = "Inactive"
$rng = $ws.Range("A43:F43")
$rng.NumberFormat = "@"
$ws.Cells.Item(43, 1).Value = "2026-01-30"
$ws.Cells.Item(43, 2).Value = "12:55:10"
$ws.Cells.Item(43, 3).Value = "12:00"
$ws.Cells.Item(43, 4).Value = "Bathroom"
$ws.Cells.Item(43, 5).Value = "No Motion"
$ws.Cells.Item(43, 6).Value = "Inactive"
$rng = $ws.Range("A44:F44")
$rng.NumberFormat = "@"
$ws.Cells.Item(44, 1).Value = "2026-01-30"
$ws.Cells.Item(44, 2).Value = "12:55:15"
$ws.Cells.Item(44, 3).Value = "12:00"
$ws.Cells.Item(44, 4).Value = "Bathroom"
$ws.Cells.Item(44, 5).Value = "No Motion"
$ws.Cells.Item(44, 6).Value = "Inactive"
$rng = $ws.Range("A45:F45")
$rng.NumberFormat = "@"
$ws.Cells.Item(45, 1).Value = "2026-01-30"
$ws.Cells.Item(45, 2).Value = "12:55:20"
$ws.Cells.Item(45, 3).Value = "12:00"
$ws.Cells.Item(45, 4).Value = "Bathroom"
$ws.Cells.Item(45, 5).Value = "No Motion"
$ws.Cells.Item(45, 6).Value = "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
$rng = $ws.Range("A27:F27")
$rng.NumberFormat = "@"
$ws.Cells.Item(27, 1).Value = "2026-01-30"
$ws.Cells.Item(27, 2).Value = "12:54:20"
$ws.Cells.Item(27, 3).Value = "12:00"
$ws.Cells.Item(27, 4).Value = "Bathroom"
$ws.Cells.Item(27, 5).Value = "87.3%"
$ws.Cells.Item(27, 6).Value = "Active"
$rng = $ws.Range("A28:F28")
$rng.NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = "2026-01-30"
$ws.Cells.Item(28, 2).Value = "12:54:21"
$ws.Cells.Item(28, 3).Value = "12:00"
$ws.Cells.Item(28, 4).Value = "Bathroom"
$ws.Cells.Item(28, 5).Value = "86.5%"
$ws.Cells.Item(28, 6).Value = "Active"
$rng = $ws.Range("A29:F29")
$rng.NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = "2026-01-30"
$ws.Cells.Item(29, 2).Value = "12:54:26"
$ws.Cells.Item(29, 3).Value = "12:00"
$ws.Cells.Item(29, 4).Value = "Bathroom"
$ws.Cells.Item(29, 5).Value = "87.4%"
$ws.Cells.Item(29, 6).Value = "Active"
$rng = $ws.Range("A30:F30")
$rng.NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "2026-01-30"
$ws.Cells.Item(30, 2).Value = "12:54:30"
$ws.Cells.Item(30, 3).Value = "12:00"
$ws.Cells.Item(30, 4).Value = "Bathroom"
$ws.Cells.Item(30, 5).Value = "87.5%"
$ws.Cells.Item(30, 6).Value = "Active"
$rng = $ws.Range("A31:F31")
$rng.NumberFormat = "@"
$ws.Cells.Item(31, 1).Value = "2026-01-30"
$ws.Cells.Item(31, 2).Value = "12:54:34"
$ws.Cells.Item(31, 3).Value = "12:00"
$ws.Cells.Item(31, 4).Value = "Bathroom"
$ws.Cells.Item(31, 5).Value = "87.5%"
$ws.Cells.Item(31, 6).Value = "Active"
$rng = $ws.Range("A32:F32")
$rng.NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "2026-01-30"
$ws.Cells.Item(32, 2).Value = "12:54:38"
$ws.Cells.Item(32, 3).Value = "12:00"
$ws.Cells.Item(32, 4).Value = "Bathroom"
$ws.Cells.Item(32, 5).Value = "86.6%"
$ws.Cells.Item(32, 6).Value = "Active"
$rng = $ws.Range("A33:F33")
$rng.NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = "2026-01-30"
$ws.Cells.Item(33, 2).Value = "12:54:46"
$ws.Cells.Item(33, 3).Value = "12:00"
$ws.Cells.Item(33, 4).Value = "Bathroom"
$ws.Cells.Item(33, 5).Value = "87.4%"
$ws.Cells.Item(33, 6).Value = "Active"
$rng = $ws.Range("A34:F34")
$rng.NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = "2026-01-30"
$ws.Cells.Item(34, 2).Value = "12:54:50"
$ws.Cells.Item(34, 3).Value = "12:00"
$ws.Cells.Item(34, 4).Value = "Bathroom"
$ws.Cells.Item(34, 5).Value = "86.6%"
$ws.Cells.Item(34, 6).Value = "Active"
$rng = $ws.Range("A35:F35")
$rng.NumberFormat = "@"
$ws.Cells.Item(35, 1).Value = "2026-01-30"
$ws.Cells.Item(35, 2).Value = "12:54:58"
$ws.Cells.Item(35, 3).Value = "12:00"
$ws.Cells.Item(35, 4).Value = "Bathroom"
$ws.Cells.Item(35, 5).Value = "86.6%"
$ws.Cells.Item(35, 6).Value = "Active"
$rng = $ws.Range("A36:F36")
$rng.NumberFormat = "@"
$ws.Cells.Item(36, 1).Value = "2026-01-30"
$ws.Cells.Item(36, 2).Value = "12:55:06"
$ws.Cells.Item(36, 3).Value = "12:00"
$ws.Cells.Item(36, 4).Value = "Bathroom"
$ws.Cells.Item(36, 5).Value = "87.5%"
$ws.Cells.Item(36, 6).Value = "Active"
$rng = $ws.Range("A37:F37")
$rng.NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = "2026-01-30"
$ws.Cells.Item(37, 2).Value = "12:55:10"
$ws.Cells.Item(37, 3).Value = "12:00"
$ws.Cells.Item(37, 4).Value = "Bathroom"
$ws.Cells.Item(37, 5).Value = "86.6%"
$ws.Cells.Item(37, 6).Value = "Active"
$rng = $ws.Range("A38:F38")
$rng.NumberFormat = "@"
$ws.Cells.Item(38, 1).Value = "2026-01-30"
$ws.Cells.Item(38, 2).Value = "12:55:14"
$ws.Cells.Item(38, 3).Value = "12:00"
$ws.Cells.Item(38, 4).Value = "Bathroom"
$ws.Cells.Item(38, 5).Value = "87.6%"
$ws.Cells.Item(38, 6).Value = "Active"
$rng = $ws.Range("A39:F39")
$rng.NumberFormat = "@"
$ws.Cells.Item(39, 1).Value = "2026-01-30"
$ws.Cells.Item(39, 2).Value = "12:55:18"
$ws.Cells.Item(39, 3).Value = "12:00"
$ws.Cells.Item(39, 4).Value = "Bathroom"
$ws.Cells.Item(39, 5).Value = "86.6%"
$ws.Cells.Item(39, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Temperature")
$rng = $ws.Range("A27:F27")
$rng.NumberFormat = "@"
$ws.Cells.Item(27, 1).Value = "2026-01-30"
$ws.Cells.Item(27, 2).Value = "12:54:20"
$ws.Cells.Item(27, 3).Value = "12:00"
$ws.Cells.Item(27, 4).Value = "Bathroom"
$ws.Cells.Item(27, 5).Value = "22.7C"
$ws.Cells.Item(27, 6).Value = "Active"
$rng = $ws.Range("A28:F28")
$rng.NumberFormat = "@"
$ws.Cells.Item(28, 1).Value = "2026-01-30"
$ws.Cells.Item(28, 2).Value = "12:54:21"
$ws.Cells.Item(28, 3).Value = "12:00"
$ws.Cells.Item(28, 4).Value = "Bathroom"
$ws.Cells.Item(28, 5).Value = "22.7C"
$ws.Cells.Item(28, 6).Value = "Active"
$rng = $ws.Range("A29:F29")
$rng.NumberFormat = "@"
$ws.Cells.Item(29, 1).Value = "2026-01-30"
$ws.Cells.Item(29, 2).Value = "12:54:26"
$ws.Cells.Item(29, 3).Value = "12:00"
$ws.Cells.Item(29, 4).Value = "Bathroom"
$ws.Cells.Item(29, 5).Value = "22.6C"
$ws.Cells.Item(29, 6).Value = "Active"
$rng = $ws.Range("A30:F30")
$rng.NumberFormat = "@"
$ws.Cells.Item(30, 1).Value = "2026-01-30"
$ws.Cells.Item(30, 2).Value = "12:54:30"
$ws.Cells.Item(30, 3).Value = "12:00"
$ws.Cells.Item(30, 4).Value = "Bathroom"
$ws.Cells.Item(30, 5).Value = "22.7C"
$ws.Cells.Item(30, 6).Value = "Active"
$rng = $ws.Range("A31:F31")
$rng.NumberFormat = "@"
$ws.Cells.Item(31, 1).Value = "2026-01-30"
$ws.Cells.Item(31, 2).Value = "12:54:34"
$ws.Cells.Item(31, 3).Value = "12:00"
$ws.Cells.Item(31, 4).Value = "Bathroom"
$ws.Cells.Item(31, 5).Value = "22.7C"
$ws.Cells.Item(31, 6).Value = "Active"
$rng = $ws.Range("A32:F32")
$rng.NumberFormat = "@"
$ws.Cells.Item(32, 1).Value = "2026-01-30"
$ws.Cells.Item(32, 2).Value = "12:54:38"
$ws.Cells.Item(32, 3).Value = "12:00"
$ws.Cells.Item(32, 4).Value = "Bathroom"
$ws.Cells.Item(32, 5).Value = "22.7C"
$ws.Cells.Item(32, 6).Value = "Active"
$rng = $ws.Range("A33:F33")
$rng.NumberFormat = "@"
$ws.Cells.Item(33, 1).Value = "2026-01-30"
$ws.Cells.Item(33, 2).Value = "12:54:46"
$ws.Cells.Item(33, 3).Value = "12:00"
$ws.Cells.Item(33, 4).Value = "Bathroom"
$ws.Cells.Item(33, 5).Value = "22.6C"
$ws.Cells.Item(33, 6).Value = "Active"
$rng = $ws.Range("A34:F34")
$rng.NumberFormat = "@"
$ws.Cells.Item(34, 1).Value = "2026-01-30"
$ws.Cells.Item(34, 2).Value = "12:54:50"
$ws.Cells.Item(34, 3).Value = "12:00"
$ws.Cells.Item(34, 4).Value = "Bathroom"
$ws.Cells.Item(34, 5).Value = "22.7C"
$ws.Cells.Item(34, 6).Value = "Active"
$rng = $ws.Range("A35:F35")
$rng.NumberFormat = "@"
$ws.Cells.Item(35, 1).Value = "2026-01-30"
$ws.Cells.Item(35, 2).Value = "12:54:58"
$ws.Cells.Item(35, 3).Value = "12:00"
$ws.Cells.Item(35, 4).Value = "Bathroom"
$ws.Cells.Item(35, 5).Value = "22.6C"
$ws.Cells.Item(35, 6).Value = "Active"
$rng = $ws.Range("A36:F36")
$rng.NumberFormat = "@"
$ws.Cells.Item(36, 1).Value = "2026-01-30"
$ws.Cells.Item(36, 2).Value = "12:55:06"
$ws.Cells.Item(36, 3).Value = "12:00"
$ws.Cells.Item(36, 4).Value = "Bathroom"
$ws.Cells.Item(36, 5).Value = "22.6C"
$ws.Cells.Item(36, 6).Value = "Active"
$rng = $ws.Range("A37:F37")
$rng.NumberFormat = "@"
$ws.Cells.Item(37, 1).Value = "2026-01-30"
$ws.Cells.Item(37, 2).Value = "12:55:10"
$ws.Cells.Item(37, 3).Value = "12:00"
$ws.Cells.Item(37, 4).Value = "Bathroom"
$ws.Cells.Item(37, 5).Value = "22.6C"
$ws.Cells.Item(37, 6).Value = "Active"
$rng = $ws.Range("A38:F38")
$rng.NumberFormat = "@"
$ws.Cells.Item(38, 1).Value = "2026-01-30"
$ws.Cells.Item(38, 2).Value = "12:55:14"
$ws.Cells.Item(38, 3).Value = "12:00"
$ws.Cells.Item(38, 4).Value = "Bathroom"
$ws.Cells.Item(38, 5).Value = "22.7C"
$ws.Cells.Item(38, 6).Value = "Active"
$rng = $ws.Range("A39:F39")
$rng.NumberFormat = "@"
$ws.Cells.Item(39, 1).Value = "2026-01-30"
$ws.Cells.Item(39, 2).Value = "12:55:18"
$ws.Cells.Item(39, 3).Value = "12:00"
$ws.Cells.Item(39, 4).Value = "Bathroom"
$ws.Cells.Item(39, 5).Value = "22.6C"
$ws.Cells.Item(39, 6).Value = "Active"

$ws = $wb.Worksheets.Item("Proximity")
$rng = $ws.Range("A8:F8")
$rng.NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "2026-01-30"
$ws.Cells.Item(8, 2).Value = "12:54:21"
$ws.Cells.Item(8, 3).Value = "12:00"
$ws.Cells.Item(8, 4).Value = "Bathroom Door"
$ws.Cells.Item(8, 5).Value = "ENTER"
$ws.Cells.Item(8, 6).Value = "User ENTERED Bathroom"
$rng = $ws.Range("A9:F9")
$rng.NumberFormat = "@"
$ws.Cells.Item(9, 1).Value = "2026-01-30"
$ws.Cells.Item(9, 2).Value = "12:54:34"
$ws.Cells.Item(9, 3).Value = "12:00"
$ws.Cells.Item(9, 4).Value = "Bathroom Door"
$ws.Cells.Item(9, 5).Value = "EXIT"
$ws.Cells.Item(9, 6).Value = "User EXITED Bathroom"
$rng = $ws.Range("A10:F10")
$rng.NumberFormat = "@"
$ws.Cells.Item(10, 1).Value = "2026-01-30"
$ws.Cells.Item(10, 2).Value = "12:54:37"
$ws.Cells.Item(10, 3).Value = "12:00"
$ws.Cells.Item(10, 4).Value = "Bathroom Door"
$ws.Cells.Item(10, 5).Value = "ENTER"
$ws.Cells.Item(10, 6).Value = "User ENTERED Bathroom"
$rng = $ws.Range("A11:F11")
$rng.NumberFormat = "@"
$ws.Cells.Item(11, 1).Value = "2026-01-30"
$ws.Cells.Item(11, 2).Value = "12:54:41"
$ws.Cells.Item(11, 3).Value = "12:00"
$ws.Cells.Item(11, 4).Value = "Bathroom Door"
$ws.Cells.Item(11, 5).Value = "EXIT"
$ws.Cells.Item(11, 6).Value = "User EXITED Bathroom"
$rng = $ws.Range("A12:F12")
$rng.NumberFormat = "@"
$ws.Cells.Item(12, 1).Value = "2026-01-30"
$ws.Cells.Item(12, 2).Value = "12:54:49"
$ws.Cells.Item(12, 3).Value = "12:00"
$ws.Cells.Item(12, 4).Value = "Bathroom Door"
$ws.Cells.Item(12, 5).Value = "ENTER"
$ws.Cells.Item(12, 6).Value = "User ENTERED Bathroom"
$rng = $ws.Range("A13:F13")
$rng.NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "2026-01-30"
$ws.Cells.Item(13, 2).Value = "12:54:53"
$ws.Cells.Item(13, 3).Value = "12:00"
$ws.Cells.Item(13, 4).Value = "Bathroom Door"
$ws.Cells.Item(13, 5).Value = "EXIT"
$ws.Cells.Item(13, 6).Value = "User EXITED Bathroom"
$rng = $ws.Range("A14:F14")
$rng.NumberFormat = "@"
$ws.Cells.Item(14, 1).Value = "2026-01-30"
$ws.Cells.Item(14, 2).Value = "12:54:56"
$ws.Cells.Item(14, 3).Value = "12:00"
$ws.Cells.Item(14, 4).Value = "Bathroom Door"
$ws.Cells.Item(14, 5).Value = "ENTER"
$ws.Cells.Item(14, 6).Value = "User ENTERED Bathroom"
$rng = $ws.Range("A15:F15")
$rng.NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "2026-01-30"
$ws.Cells.Item(15, 2).Value = "12:55:03"
$ws.Cells.Item(15, 3).Value = "12:00"
$ws.Cells.Item(15, 4).Value = "Bathroom Door"
$ws.Cells.Item(15, 5).Value = "EXIT"
$ws.Cells.Item(15, 6).Value = "User EXITED Bathroom"
$rng = $ws.Range("A16:F16")
$rng.NumberFormat = "@"
$ws.Cells.Item(16, 1).Value = "2026-01-30"
$ws.Cells.Item(16, 2).Value = "12:55:13"
$ws.Cells.Item(16, 3).Value = "12:00"
$ws.Cells.Item(16, 4).Value = "Bathroom Door"
$ws.Cells.Item(16, 5).Value = "ENTER"
$ws.Cells.Item(16, 6).Value = "User ENTERED Bathroom"
$rng = $ws.Range("A17:F17")
$rng.NumberFormat = "@"
$ws.Cells.Item(17, 1).Value = "2026-01-30"
$ws.Cells.Item(17, 2).Value = "12:55:16"
$ws.Cells.Item(17, 3).Value = "12:00"
$ws.Cells.Item(17, 4).Value = "Bathroom Door"
$ws.Cells.Item(17, 5).Value = "EXIT"
$ws.Cells.Item(17, 6).Value = "User EXITED Bathroom"

$ws = $wb.Worksheets.Item("mmWave")
$rng = $ws.Range("A19:F19")
$rng.NumberFormat = "@"
$ws.Cells.Item(19, 1).Value = "2026-01-30"
$ws.Cells.Item(19, 2).Value = "12:54:33"
$ws.Cells.Item(19, 3).Value = "12:00"
$ws.Cells.Item(19, 4).Value = "Living Room"
$ws.Cells.Item(19, 5).Value = "PRESENCE_DETECTED"
$ws.Cells.Item(19, 6).Value = "Active"
$rng = $ws.Range("A20:F20")
$rng.NumberFormat = "@"
$ws.Cells.Item(20, 1).Value = "2026-01-30"
$ws.Cells.Item(20, 2).Value = "12:54:57"
$ws.Cells.Item(20, 3).Value = "12:00"
$ws.Cells.Item(20, 4).Value = "Living Room"
$ws.Cells.Item(20, 5).Value = "FALL_DETECTED"
$ws.Cells.Item(20, 6).Value = "EMERGENCY"
